$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 3 (proto-refseqs flanking data) - set first so its new shared string
# (ENSHCOG00000009382) lands at the lowest new index, matching the target order.
$ws.Range("P3").Value = "ENSHCOG00000009382"
$ws.Range("R3").Value = 16

# Row 2 (ichthama-refseqs flanking data)
$ws.Range("P2").Value = "BHLHE23"
$ws.Range("Q2").Value = "YTHDF1"
$ws.Range("R2").Value = 20
